# Commit: Tue, May 12, 2020 10:05:46 PM
#
# The deck's applied design theme ("Integral") is switched back to the
# default Office theme. The slide master (and therefore every slide /
# layout, since this deck has a single master) picks up the stock
# "Office" theme color scheme:
#
#   dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink
#     000000 / FFFFFF / 44546A / E7E6E6 / 5B9BD5 / ED7D31 / A5A5A5 /
#     FFC000 / 4472C4 / 70AD47 / 0563C1 / 954F72
#
# (font scheme and format scheme - fills/lines/effects - are identical
# between "Integral" and the default "Office" theme, so only the 12
# theme colors actually move.)
#
# Walking the theme color scheme off a Slide (rather than the Master's
# old-style 8-color ColorScheme) updates every one of the twelve
# DrawingML theme colors in place on the deck's live theme part without
# disturbing anything else.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      -> #000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> #FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> #44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> #E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> #5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> #ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> #A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> #FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> #4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> #70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> #0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> #954F72
